$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused | Clear Glass Lens
$ws.Range("H33").Value = 3034.162
$ws.Range("I33").Value = 52.666668
$ws.Range("J33").Value = 5067
$ws.Range("K33").Value = 52.666668
$ws.Range("L33").Value = 5067
$ws.Range("M33").Value = 176.333332
$ws.Range("N33").Value = -5525

# Row 100: Asking for a Friend | Beetle Glue
$ws.Range("H100").Value = 15153511
$ws.Range("I100").Value = 23810612
$ws.Range("J100").Value = 3585.75
$ws.Range("K100").Value = 23810612
$ws.Range("L100").Value = 3585.75
$ws.Range("M100").Value = -23810071
$ws.Range("N100").Value = -4667.75

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 1117.8334
$ws.Range("I137").Value = 901.61536
$ws.Range("J137").Value = 1680
$ws.Range("K137").Value = 2704.84608
$ws.Range("L137").Value = 5040
$ws.Range("M137").Value = -154.8460800000003
$ws.Range("N137").Value = -10140

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 4786
$ws.Range("I141").Value = 4849.2856
$ws.Range("K141").Value = 14547.8568
$ws.Range("M141").Value = -9367.856800000001

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 964.875
$ws.Range("I74").Value = 834.3043
$ws.Range("J74").Value = 1141.5294
$ws.Range("K74").Value = 834.3043
$ws.Range("L74").Value = 1141.5294
$ws.Range("M74").Value = 39.69569999999999
$ws.Range("N74").Value = -2889.5294

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 964.875
$ws.Range("I77").Value = 834.3043
$ws.Range("J77").Value = 1141.5294
$ws.Range("K77").Value = 4171.5215
$ws.Range("L77").Value = 5707.646999999999
$ws.Range("M77").Value = 196.4785000000002
$ws.Range("N77").Value = -14443.647

$ws = $wb.Worksheets.Item("BSM")
# Row 98: Killer Cutlery | Doman Iron Culinary Knife
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 6979.0527
$ws.Range("I31").Value = 7211.8823
$ws.Range("K31").Value = 7211.8823
$ws.Range("M31").Value = -6916.8823

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 6979.0527
$ws.Range("I34").Value = 7211.8823
$ws.Range("K34").Value = 7211.8823
$ws.Range("M34").Value = -7009.8823

# Row 140: Spear Pressure | Claro Walnut Spear
$ws.Range("H140").Value = 58666.668
$ws.Range("I140").Value = 10000
$ws.Range("J140").Value = 83000
$ws.Range("K140").Value = 10000
$ws.Range("L140").Value = 83000
$ws.Range("M140").Value = -4820
$ws.Range("N140").Value = -93360

$ws = $wb.Worksheets.Item("CUL")
# Row 110: His Dark Utensils | Spaghetti al Nero
$ws.Range("H110").Value = 5200
$ws.Range("I110").Value = 5200
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 15600
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -11510
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 42: It's Only Love | Silver Choker
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 99: Needle in a Hingan Stack | Dzo Horn Needle
$ws.Range("H99").Value = 3300.8125
$ws.Range("I99").Value = 3300.8125
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3300.8125
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1054.8125
$ws.Range("N99").ClearContents()

# Row 104: Speak Softly and Carry a Metal Rod | Palladium Rod
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 105: Untucked | Palladium Tuck
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 107: Whetstones for the Workers | Hard Mudstone Whetstone
$ws.Range("H107").Value = 346.45456
$ws.Range("I107").Value = 346.45456
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 346.45456
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1573.54544
$ws.Range("N107").ClearContents()

# Row 108: Satisfactory Sewing | Stonegold Needle
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 109: You're My Wonderhall | Hematite Earrings of Healing
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 110: Slimming Down | Stonegold Rapier
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 1414.45
$ws.Range("I113").Value = 1236.8125
$ws.Range("J113").Value = 2125
$ws.Range("K113").Value = 1236.8125
$ws.Range("L113").Value = 2125
$ws.Range("M113").Value = 933.1875
$ws.Range("N113").Value = -6465

# Row 114: Hot Rod | Bluespirit Rod
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 115: Unsung Generosity | Manasilver Choker
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# Row 117: Birth Ring | Triplite Ring of Aiming
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 118: A Magnanimous Refrain | Triplite Earrings of Casting
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# Row 121: Wrap Those Wrists | Petalite Bracelet of Fending
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 1243.0454
$ws.Range("I122").Value = 1207
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3621
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1171
$ws.Range("N122").Value = -10900

# Row 123: Workplace Workout | Ametrine Ring of Fending
$ws.Range("H123").Value = 14119.546
$ws.Range("J123").Value = 14119.546
$ws.Range("L123").Value = 14119.546
$ws.Range("N123").Value = -19019.546

# Row 124: The Sage's Successor | Pewter Pendulums
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Row 125: Pewter-hewn Punishment | Pewter Choker of Slaying
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 11118661
$ws.Range("I126").Value = 9849.091
$ws.Range("J126").Value = 41667892
$ws.Range("K126").Value = 29547.273
$ws.Range("L126").Value = 125003676
$ws.Range("M126").Value = -27077.273
$ws.Range("N126").Value = -125008616

# Row 127: Sage with the Golden Earrings | Phrygian Ear Cuffs of Healing
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# Row 128: To Fight at Her Side | Manganese Rapier
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 130: Planisphere to Paper | Chondrite Magitek Planisphere
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 69151.336
$ws.Range("I132").Value = 1807.25
$ws.Range("J132").Value = 146116
$ws.Range("K132").Value = 5421.75
$ws.Range("L132").Value = 438348
$ws.Range("M132").Value = -2891.75
$ws.Range("N132").Value = -443408

# Row 133: Pendulums of Our Own | Lar Pendulums
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 137: Sew Excited | Cobalt Tungsten Needle
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 50: The Birdmen of Ishgard | Boarskin Culottes
$ws.Range("H50").Value = 5500
$ws.Range("J50").Value = 7000
$ws.Range("L50").Value = 7000
$ws.Range("N50").Value = -8274

# Row 62: Pummeling Abroad | Archaeoskin Breeches of Maiming
$ws.Range("H62").Value = 10933.333
$ws.Range("J62").Value = 10933.333
$ws.Range("L62").Value = 10933.333
$ws.Range("N62").Value = -12181.333

# Row 63: From Mud to Mourning | Archaeoskin Jackboots of Gathering
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 65: The Style of the Time (L) | Archaeoskin Breeches of Maiming
$ws.Range("H65").Value = 10933.333
$ws.Range("J65").Value = 10933.333
$ws.Range("L65").Value = 32799.999
$ws.Range("N65").Value = -39039.999

# Row 66: These Boots Are Made for Hawkin' (L) | Archaeoskin Jackboots of Gathering
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 74: Overall, We Blend In | Dhalmelskin Vest
$ws.Range("H74").Value = 15850
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 15850
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 15850
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -17846

# Row 77: Eviction Notice (L) | Dhalmelskin Vest
$ws.Range("H77").Value = 15850
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 15850
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 47550
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -57534

# Row 88: It Will Knock Your Socks Off | Serpentskin Thighboots of Casting
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

# Row 91: On My Own Two Feet (L) | Serpentskin Thighboots of Casting
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 18121.334
$ws.Range("I136").Value = 21880.8
$ws.Range("J136").Value = 13422
$ws.Range("K136").Value = 65642.39999999999
$ws.Range("L136").Value = 40266
$ws.Range("M136").Value = -63092.39999999999
$ws.Range("N136").Value = -45366
